$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column E ("Tong tien qua tang") to account for
# gift-product cost in the branch revenue export, shifting the remaining
# revenue columns (old E..O) one slot to the right (new F..P).
$ws.Columns("E:E").Insert()
$ws.Columns("E:E").ColumnWidth = $ws.Columns("D:D").ColumnWidth

# New header label for the inserted column.
$ws.Range("E8").Value = "Tổng tiền quà tặng"

# Restore the user's cell selection as last saved in the source workbook.
$ws.Range("J15").Select()
